$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entries continuing the effort log (manual chapter API drafted)
$ws.Range("A28").Value = 41200
$ws.Range("B28").Value = 1.75
$ws.Range("D28").Value = "Manual continued"

$ws.Range("A29").Value = 41205
$ws.Range("B29").Value = 2.5
$ws.Range("D29").Value = "Manual continued"

# Update view: scroll so the new rows are visible, select the last new cell
$ws.Range("D29").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
